# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md file (row 7 of
# each sheet) to reflect a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-19 04:37:05"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-19 04:36:57"

# de-de sheet: column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-19 04:37:05"
